# Updates cryptos list price (D) and volume-change (E) columns for rows 2-51
# per the scraped GitHub Actions refresh run.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.100.05"
$ws.Range("E2").Value = "  -1.39%  "

$ws.Range("D3").Value = "1.897.83"
$ws.Range("E3").Value = "  -0.55%  "

$ws.Range("D4").Value = "'1.000"
$ws.Range("E4").Value = "  -0.06%  "

$ws.Range("D5").Value = "'314.26"
$ws.Range("E5").Value = "  +0.17%  "

$ws.Range("D6").Value = "'0.9999"
$ws.Range("E6").Value = "  -0.05%  "

$ws.Range("D7").Value = "'0.5031"
$ws.Range("E7").Value = "  -0.38%  "

$ws.Range("D8").Value = "'0.3899"
$ws.Range("E8").Value = "  -1.30%  "

$ws.Range("D9").Value = "'0.09251"
$ws.Range("E9").Value = "  -5.56%  "

$ws.Range("E10").Value = "  -2.76%  "

$ws.Range("D11").Value = "'41.78"
$ws.Range("E11").Value = "  +0.56%  "

$ws.Range("D12").Value = "'6.405"
$ws.Range("E12").Value = "  -2.20%  "

$ws.Range("D13").Value = "'20.82"
$ws.Range("E13").Value = "  -1.49%  "

$ws.Range("D14").Value = "1.891.15"
$ws.Range("E14").Value = "  -1.06%  "

$ws.Range("D15").Value = "'7.291"
$ws.Range("E15").Value = "  -3.78%  "

$ws.Range("D16").Value = "'0.9997"
$ws.Range("E16").Value = "  -0.09%  "

$ws.Range("D17").Value = "'92.44"
$ws.Range("E17").Value = "  -1.50%  "

$ws.Range("D18").Value = "'0.00001109"
$ws.Range("E18").Value = "  -2.95%  "

$ws.Range("D19").Value = "'0.06640"
$ws.Range("E19").Value = "  -0.23%  "

$ws.Range("D20").Value = "'17.86"
$ws.Range("E20").Value = "  -1.26%  "

$ws.Range("E21").Value = "  -0.05%  "

$ws.Range("D22").Value = "'6.225"
$ws.Range("E22").Value = "  -1.19%  "

$ws.Range("D23").Value = "28.156.23"
$ws.Range("E23").Value = "  -1.37%  "

$ws.Range("D24").Value = "'11.45"
$ws.Range("E24").Value = "  +0.06%  "

$ws.Range("D25").Value = "'2.310"
$ws.Range("E25").Value = "  +1.27%  "

$ws.Range("D26").Value = "2.112.72"
$ws.Range("E26").Value = "  -0.74%  "

$ws.Range("D27").Value = "'2.557"
$ws.Range("E27").Value = "  -6.68%  "

$ws.Range("D28").Value = "'20.88"
$ws.Range("E28").Value = "  -2.25%  "

$ws.Range("D29").Value = "'158.13"
$ws.Range("E29").Value = "  -0.78%  "

$ws.Range("D30").Value = "'126.49"
$ws.Range("E30").Value = "  -1.85%  "

$ws.Range("E31").Value = "  -0.95%  "

$ws.Range("D32").Value = "'0.1062"
$ws.Range("E32").Value = "  -0.86%  "

$ws.Range("D33").Value = "'5.617"
$ws.Range("E33").Value = "  -1.69%  "

$ws.Range("D34").Value = "'3.614"
$ws.Range("E34").Value = "  -0.59%  "

$ws.Range("D35").Value = "'9.604"
$ws.Range("E35").Value = "  -3.37%  "

$ws.Range("D36").Value = "'0.06604"
$ws.Range("E36").Value = "  -2.78%  "

$ws.Range("D37").Value = "'0.02402"
$ws.Range("E37").Value = "  -1.78%  "

$ws.Range("D38").Value = "'0.2207"
$ws.Range("E38").Value = "  -1.24%  "

$ws.Range("D39").Value = "'1.301"
$ws.Range("E39").Value = "  +9.03%  "

$ws.Range("D40").Value = "'1.226"
$ws.Range("E40").Value = "  -3.81%  "

$ws.Range("D41").Value = "'0.6489"
$ws.Range("E41").Value = "  +0.65%  "

$ws.Range("D42").Value = "'4.980"
$ws.Range("E42").Value = "  -2.28%  "

$ws.Range("D43").Value = "'11.42"
$ws.Range("E43").Value = "  -1.92%  "

$ws.Range("D44").Value = "'0.9996"
$ws.Range("E44").Value = "  -0.06%  "

$ws.Range("D45").Value = "'0.6109"
$ws.Range("E45").Value = "  -0.13%  "

$ws.Range("D46").Value = "'13.29"
$ws.Range("E46").Value = "  -2.79%  "

$ws.Range("E47").Value = "  +2.04%  "

$ws.Range("D48").Value = "'3.692"
$ws.Range("E48").Value = "  +0.82%  "

$ws.Range("D49").Value = "'2.003"
$ws.Range("E49").Value = "  -2.20%  "

$ws.Range("D50").Value = "'122.32"
$ws.Range("E50").Value = "  -2.49%  "

$ws.Range("D51").Value = "'1.189"
$ws.Range("E51").Value = "  -1.80%  "
